$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Change 1: "NRC Postdoc," -> "Postdoc Associate," (bold run), and
# collapse the long run of spacer runs between "Boulder, CO" and
# "Feb. 2017- 2019" into a tab + a block of small (sz4) spaces.
# -------------------------------------------------------------------

# 1a. Rename the bold title text.
$r1 = $d.Content
$r1.Find.Execute("NRC Postdoc,", $true, $false, $false, $false, $false, $true, 1, $false, "Postdoc Associate,", 2)

# 1b. Replace the whitespace stretch between "Boulder, CO" and "Feb."
#     with a single tab character followed by 50 narrow (sz4) spaces.
$tab = [char]9
$spacer = $tab + "                                                  "
$r2 = $d.Content
$r2.Find.Execute("Boulder, CO[ ]{1,}Feb", $true, $false, $true, $false, $false, $true, 1, $false, ("Boulder, CO" + $spacer + "Feb"), 2)

# -------------------------------------------------------------------
# Change 2: split the NRC sentence to parenthesize "(NRC)".
# -------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute("National Research Council selection committee", $true, $false, $false, $false, $false, $true, 1, $false, "National Research Council (NRC) selection committee", 2)

# -------------------------------------------------------------------
# Change 3: the _GoBack bookmark moved from the "intermediate" skill
# bullet (GPU Optimized Monte Carlo) to the NRC scoring bullet above.
# Remove it from its old location ...
# -------------------------------------------------------------------
foreach ($bm in @($d.Bookmarks)) {
    if ($bm.Name -eq "_GoBack") {
        $bm.Delete()
    }
}

# ... and add it back at the new location (end of the "(NRC)" bullet,
# right before " selection committee").
$r4 = $d.Content
$r4.Find.Execute("(NRC)", $false)
$bmRange = $d.Range($r4.End, $r4.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
